# ARKCORR-22 Minor changes to the dueDate setting for the 'Release' queue.
#
# The "Set Due Date Release Queue" rule (row 31) in Sheet1 is updated so the
# CONDITION no longer checks "dueDate == null" (only the queue name matters
# now) and the ACTION simply clears the due date instead of computing
# today's date.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# CONDITION column (C) for the Release-queue rule
$ws.Range("C31").Value = 'queue.name == "Release"'

# ACTION column (D) for the Release-queue rule
$ws.Range("D31").Value = "setDueDate, null"

# Match the author's final cursor position (previously on C31)
$ws.Range("D31").Select()
